$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (H1:K1)
$ws.Range("H1").Value = "Labor Booking User"
$ws.Range("I1").Value = "SiteID"
$ws.Range("J1").Value = "Location ID"
$ws.Range("K1").Value = "Location Number"

# New data cells for row 2
$ws.Range("H2").Value = "a811K0000004fpN"
$ws.Range("I2").Value = "a7q410000004I1W"
$ws.Range("J2").Value = "a7Z4100000000hb"
$ws.Range("K2").Value = "SY_ReceiptLoc"

# New data cells for row 3
$ws.Range("H3").Value = "a811K0000004fpN"
$ws.Range("I3").Value = "a7q410000004I1W"
$ws.Range("J3").Value = "a7Z4100000000hb"
$ws.Range("K3").Value = "SY_ReceiptLoc"

# Header row A1:F1 loses its bold formatting (reverts to default style)
$ws.Range("A1:F1").ClearFormats()

# Column width adjustments
$ws.Columns("B").ColumnWidth = 21.44140625
$ws.Columns("H").ColumnWidth = 16.109375
$ws.Columns("I").ColumnWidth = 16.44140625
$ws.Columns("J").ColumnWidth = 16
$ws.Columns("K").ColumnWidth = 12.88671875

# Selection moves to I2
$ws.Range("I2").Select()
